$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current last data row (row 45), which
# pushes the existing row 45 ("27", "Partly Cloudy", "01/18/2025", "20")
# down to row 48 untouched.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# Fill the three newly inserted rows (45-47) with the same data as the
# (now relocated) last row, but store the temperature and hour as real
# numbers instead of text. The date is entered with a leading apostrophe
# so Excel keeps it as plain text ("01/18/2025") instead of auto-
# converting it into a date serial value.
for ($r = 45; $r -le 47; $r++) {
    $ws.Cells.Item($r, 1).Value = 27
    $ws.Cells.Item($r, 2).Value = "Partly Cloudy"
    $ws.Cells.Item($r, 3).Value = "'01/18/2025"
    $ws.Cells.Item($r, 4).Value = 20
}

Write-Host "done"
